$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet/tab
$ws.Name = "AlphaFiberF"

# 2. Fix a floating point rounding value in H15
$ws.Range("H15").Value = 0.997957694070406

# 3. Add new row 16 with data for index 14 / HexGrid-60degTilt5degRes
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9971534375439213
$ws.Range("D16").Value = 1.00908493527576
$ws.Range("E16").Value = 0.9901014836435641
$ws.Range("F16").Value = 0.9971534375439213
$ws.Range("G16").Value = 1.001074986944671
$ws.Range("H16").Value = 0.9811833402172897
$ws.Range("I16").Value = 0.9918300717683775
$ws.Range("J16").Value = 1.00908493527576
$ws.Range("K16").Value = 0.9995932094596618
$ws.Range("L16").Value = 0.9983733235017915
$ws.Range("M16").Value = 0.9950713758989305

# Match style of A16 to A15 (bold/centered/bordered style used for column A)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
